# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.225.05"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.928.31"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'248.64"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'0.7115"
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'0.3206"
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("D9").Value = "'27.25"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "'0.07078"
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("D11").Value = "'0.7901"
$ws.Range("E11").Value = "  -2.16%  "
$ws.Range("D12").Value = "'0.07986"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").Value = "1.930.46"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "'5.362"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "'94.83"
$ws.Range("D16").Value = "'14.62"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "30.245.02"
$ws.Range("D18").Value = "'254.71"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "'0.000008013"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'5.744"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").Value = "2.183.74"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'6.802"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").Value = "'9.521"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("D26").Value = "'166.03"
$ws.Range("E26").Value = "  +4.14%  "
$ws.Range("D27").Value = "'19.03"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").Value = "'2.269"
$ws.Range("E28").Value = "  -5.30%  "
$ws.Range("D29").Value = "'0.1269"
$ws.Range("E29").Value = "  -5.53%  "
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").Value = "'1.527"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").Value = "'4.383"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "'4.113"
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("D34").Value = "'0.05145"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").Value = "'1.262"
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("D36").Value = "'0.7425"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").Value = "'2.768"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").Value = "'0.01945"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").Value = "'2.802"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").Value = "'77.58"
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("D41").Value = "'6.316"
$ws.Range("E41").Value = "  -4.50%  "
$ws.Range("D42").Value = "'0.4468"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'1.969"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").Value = "'0.8452"
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "'100.47"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").Value = "'9.677"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "'7.415"
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("D49").Value = "'36.30"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Value = "'0.06119"
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'2.863"
$ws.Range("E51").Value = "  +8.75%  "
